# Api Queries for PFM Added. (#33)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 23: PFM / Transaction Details / SCPAccounts / GET / endpoint / issue
$ws.Range("A23").Value = "PFM"
$ws.Range("B23").Value = "Transaction Details"
$ws.Range("C23").Value = "SCPAccounts"
$ws.Range("D23").Value = "GET"
$ws.Range("E23").Value = "/scp/account/transactions "
$ws.Range("F23").Value = "ChildTransactiosns are what   basis for the API:- account/transactions in SCPAccounts.yaml"

# Additional rows with just the Stream column filled in
$ws.Range("A24").Value = "PFM"
$ws.Range("A25").Value = "PFM"

# Update the visible selection to match the new editing location
$ws.Range("C25").Select()
